# Insert a new "category" column before the current "author" column (D),
# shifting author/publisher/file/keywords/note/abstract one column to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Insert()

# New column D width (per target layout). The COM ColumnWidth setter in
# this runtime adds a constant ~0.8333 offset to the stored XML width, so
# subtract it here to land exactly on the target "width" attribute value.
$ws.Columns("D").ColumnWidth = 18.166666666666668
# Column I (old H = "note") width changed slightly in the target layout
$ws.Columns("I").ColumnWidth = 21.166666666666668

# Header for the new column
$ws.Range("D1").Value = "category"

# Populate the new "category" column for each paper row
$ws.Range("D2").Value = "VDU"
$ws.Range("D3").Value = "Transformer Language Model"
$ws.Range("D4").Value = "VDU"
$ws.Range("D5").Value = "Transformer Language Model"
$ws.Range("D6").Value = "Transformer Language Model"
$ws.Range("D7").Value = "VDU"
$ws.Range("D8").Value = "VDU"
$ws.Range("D9").Value = "VDU"

# The "note" column for the DocFormer row (now column I after the shift) gets new text
$ws.Range("I2").Value = "Good image for comparison of VDU methods on page 2"
